# "Fruta / hortaliza, semanal" — weekly refresh of the Pepino ensalada
# (cucumber) dataset: a new price-survey record is inserted at row 416,
# pushing the existing rows 416-463 down to 417-464 (dimension grows from
# A1:R463 to A1:R464).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 416..463 down to 417..464, leaving row 416 free for the new
# weekly record (Excel's native row-insert semantics, like right-clicking
# row 416 and choosing "Insert").
$ws.Rows.Item(416).Insert()

# Populate the newly inserted row 416 with this week's record.
$ws.Range("A416").Value = 3
$ws.Range("B416").Value = "Femacal de La Calera"
$ws.Range("C416").Value = "Coquimbo"
$ws.Range("D416").Value2 = 44918
$ws.Range("E416").Value = 5
$ws.Range("F416").Value = 100112043
$ws.Range("G416").Value = "Pepino ensalada"
$ws.Range("H416").Value = "Sin especificar"
$ws.Range("I416").Value = "Primera"
$ws.Range("J416").Value = 123
$ws.Range("K416").Value = 14500
$ws.Range("L416").Value = 15000
$ws.Range("M416").Value = 14764
$ws.Range("N416").Value = "$/caja 60 unidades"
$ws.Range("O416").Value = "Limache"
$ws.Range("P416").Value = 246
$ws.Range("Q416").Value = 60
$ws.Range("R416").Value = "Hortaliza"
